# "Working on Receipt TAB"
# Update the Application No value on the Commercial sheet and move the
# active selection, matching the author's in-progress edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commercial")
$ws.Activate()

# Application No cell changed from JP30000166 to JP30000177
$ws.Range("AA2").Value = "JP30000177"

# Active selection moved from Y9 to AF16 while working on the tab
$ws.Range("AF16").Select()
